# Auto-generated Excel COM-interop edit script
# Refactor flashcard content: replace Q/A text with revised wording,
# trim the deck to 12 Q/A rows (clear rows 14-23), fix bold styling on
# rows 7-8, and update the sheet view selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update question/answer text + row heights for rows 2-13 ---
$ws.Range("A2").Value2 = "Define Information Technology (IT)."
$ws.Range("B2").Value2 = "Information Technology (IT) is the use of computers to make humans more productive by allowing them to create, store, process and retrieve information. It encompasses hardware, software, networks and services that enable these functions."
$ws.Rows.Item(2).RowHeight = 72

$ws.Range("A3").Value2 = "Explain how Information Technology (IT) enhances human productivity in everyday life."
$ws.Range("B3").Value2 = "- It automates routine tasks such as spreadsheets for budgets - It provides instant access to information through email and search engines - It enables collaboration via cloud document sharing - It supports creativity with digital design tools"
$ws.Rows.Item(3).RowHeight = 72

$ws.Range("A4").Value2 = "Compare and contrast Information Technology (IT) and Computer Science (CS)."
$ws.Range("B4").Value2 = "- Information Technology (IT) focuses on applying and managing computing technologies to solve practical problems such as installing networks and configuring systems - Computer Science (CS) focuses on the theory and design of algorithms, data structures and software, for example writing new programming languages and researching artificial intelligence"
$ws.Rows.Item(4).RowHeight = 100.8

$ws.Range("A5").Value2 = "Define Computer Science (CS)."
$ws.Range("B5").Value2 = "Computer Science (CS) is the study of algorithms, data structures, programming languages, artificial intelligence and machine learning, as well as the theoretical foundations of computing hardware and software."
$ws.Rows.Item(5).RowHeight = 72

$ws.Range("A6").Value2 = "Define Software Engineering (SE)."
$ws.Range("B6").Value2 = "Software Engineering (SE) is the disciplined, engineering-based approach to designing, developing, testing and maintaining software applications using principles of engineering, best-practice design and programming languages to build reliable, scalable software for end users."
$ws.Rows.Item(6).RowHeight = 86.4

$ws.Range("A7").Value2 = "Compare Computer Science (CS) versus Software Engineering (SE)."
$ws.Range("B7").Value2 = "- Computer Science (CS) emphasizes theoretical foundations such as algorithms and computation theory - Software Engineering (SE) emphasizes the practical application of engineering principles to build and maintain large software systems, including software development lifecycle, testing and project management"
$ws.Rows.Item(7).RowHeight = 86.4

$ws.Range("A8").Value2 = "Define Information System (IS) and list its main components."
$ws.Range("B8").Value2 = "Information System (IS) is a collection of hardware, software, data, people and networks designed to collect, process, store and disseminate information for a specific purpose. - Hardware: PCs, servers, routers - Software: applications, operating systems - Data: raw facts and processed information - People: users and IT staff - Networks: communication links"
$ws.Rows.Item(8).RowHeight = 100.8

$ws.Range("A9").Value2 = "List and briefly describe three common types of Information Systems (IS)."
$ws.Range("B9").Value2 = "- Transaction Processing System (TPS): automates routine business transactions such as point-of-sale - Management Information System (MIS): provides reports and dashboards for decision making such as monthly sales summaries - Decision Support System (DSS): offers analytical tools and models to support complex decisions such as what-if scenario analysis"
$ws.Rows.Item(9).RowHeight = 100.8

$ws.Range("A10").Value2 = "Define Cognitive Science and name three disciplines that inform it."
$ws.Range("B10").Value2 = "Cognitive Science is the study of how the human brain perceives, learns and stores information. It draws on Psychology (mental processes), Neuroscience (brain function) and Linguistics (language processing)."
$ws.Rows.Item(10).RowHeight = 57.6

$ws.Range("A11").Value2 = "What is Computer Engineering and what are typical specializations?"
$ws.Range("B11").Value2 = "Computer Engineering blends electrical engineering and computer science to design and implement computer hardware and low-level software. Typical specializations include embedded systems for dedicated devices, very-large-scale integration (VLSI) for designing integrated circuits and networking for designing local and wide area communication systems."
$ws.Rows.Item(11).RowHeight = 100.8

$ws.Range("A12").Value2 = "Explain the measures of central tendency and when to use each."
$ws.Range("B12").Value2 = "- Mean: sum of values divided by number of values, used for symmetric distributions - Median: middle value when data are ordered, robust against outliers, used when data are skewed - Mode: most frequent value, used for categorical data or multimodal distributions"
$ws.Rows.Item(12).RowHeight = 72

$ws.Range("A13").Value2 = "Describe the differences between bar graphs, line graphs and pie charts."
$ws.Range("B13").Value2 = "- Bar graph: uses vertical or horizontal bars to compare quantities across categories - Line graph: connects data points to show trends over time or continuous variables - Pie chart: divides a circle into slices representing proportions of a whole, used to show percentage breakdowns of a single variable"
$ws.Rows.Item(13).RowHeight = 86.4

# --- Rows 7 and 8 answers are no longer bold (style matches the rest) ---
$ws.Range("B7").Font.Bold = $false
$ws.Range("B8").Font.Bold = $false

# --- The deck now only has 12 Q/A pairs: clear out the old rows 14-23 ---
$ws.Range("A14:B23").ClearContents()
$ws.Range("A14:B23").EntireRow.AutoFit()

# --- Update the saved view/selection to match the new scroll position ---
$ws.Activate()
$ws.Range("E4").Select() | Out-Null
try {
    $excel.ActiveWindow.ScrollRow = 11
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # older/limited hosts may not expose window scroll position - non-fatal
}

Write-Host "Flashcard refactor applied."
